$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) "29.565.59"
Set-TextValue $ws.Cells.Item(2, 5) "  -0.44%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.852.14"
Set-TextValue $ws.Cells.Item(3, 5) "  -0.08%  "
Set-TextValue $ws.Cells.Item(4, 4) "0.9994"
Set-TextValue $ws.Cells.Item(5, 4) "243.04"
Set-TextValue $ws.Cells.Item(5, 5) "  -0.39%  "
Set-TextValue $ws.Cells.Item(6, 4) "0.6368"
Set-TextValue $ws.Cells.Item(6, 5) "  -0.31%  "
Set-TextValue $ws.Cells.Item(8, 2) "OKB"
Set-TextValue $ws.Cells.Item(8, 3) "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Cells.Item(8, 4) "48.37"
Set-TextValue $ws.Cells.Item(8, 5) "  +3.17%  "
Set-TextValue $ws.Cells.Item(9, 2) "Dogecoin"
Set-TextValue $ws.Cells.Item(9, 3) "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Cells.Item(9, 4) "0.07592"
Set-TextValue $ws.Cells.Item(9, 5) "  +1.61%  "
Set-TextValue $ws.Cells.Item(10, 2) "Cardano"
Set-TextValue $ws.Cells.Item(10, 3) "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Cells.Item(10, 4) "0.3008"
Set-TextValue $ws.Cells.Item(10, 5) "  +0.94%  "
Set-TextValue $ws.Cells.Item(11, 2) "Solana"
Set-TextValue $ws.Cells.Item(11, 3) "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Cells.Item(11, 4) "24.23"
Set-TextValue $ws.Cells.Item(11, 5) "  -0.39%  "
Set-TextValue $ws.Cells.Item(12, 2) "TRON"
Set-TextValue $ws.Cells.Item(12, 3) "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Cells.Item(12, 4) "0.07686"
Set-TextValue $ws.Cells.Item(12, 5) "  +0.60%  "
Set-TextValue $ws.Cells.Item(13, 2) "WrappedEther"
Set-TextValue $ws.Cells.Item(13, 3) "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Cells.Item(13, 4) "1.855.20"
Set-TextValue $ws.Cells.Item(13, 5) "  +0.69%  "
Set-TextValue $ws.Cells.Item(14, 2) "Polkadot"
Set-TextValue $ws.Cells.Item(14, 3) "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Cells.Item(14, 4) "5.043"
Set-TextValue $ws.Cells.Item(14, 5) "  +0.06%  "
Set-TextValue $ws.Cells.Item(15, 2) "Polygon"
Set-TextValue $ws.Cells.Item(15, 3) "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Cells.Item(15, 4) "0.6885"
Set-TextValue $ws.Cells.Item(15, 5) "  +0.11%  "
Set-TextValue $ws.Cells.Item(16, 2) "Litecoin"
Set-TextValue $ws.Cells.Item(16, 3) "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Cells.Item(16, 4) "84.13"
Set-TextValue $ws.Cells.Item(16, 5) "  +0.49%  "
Set-TextValue $ws.Cells.Item(17, 2) "ShibaInu"
Set-TextValue $ws.Cells.Item(17, 3) "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Cells.Item(17, 4) "0.000009772"
Set-TextValue $ws.Cells.Item(17, 5) "  +3.13%  "
Set-TextValue $ws.Cells.Item(18, 2) "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Cells.Item(18, 3) "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Cells.Item(18, 4) "2.098.15"
Set-TextValue $ws.Cells.Item(18, 5) "  -0.13%  "
Set-TextValue $ws.Cells.Item(19, 2) "Uniswap"
Set-TextValue $ws.Cells.Item(19, 3) "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Cells.Item(19, 4) "6.315"
Set-TextValue $ws.Cells.Item(19, 5) "  +4.43%  "
Set-TextValue $ws.Cells.Item(20, 2) "WrappedBTC"
Set-TextValue $ws.Cells.Item(20, 3) "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Cells.Item(20, 4) "29.570.72"
Set-TextValue $ws.Cells.Item(20, 5) "  -0.47%  "
Set-TextValue $ws.Cells.Item(21, 2) "BitcoinCash"
Set-TextValue $ws.Cells.Item(21, 3) "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Cells.Item(21, 4) "238.49"
Set-TextValue $ws.Cells.Item(21, 5) "  +1.55%  "
Set-TextValue $ws.Cells.Item(22, 2) "Avalanche"
Set-TextValue $ws.Cells.Item(22, 3) "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Cells.Item(22, 4) "12.56"
Set-TextValue $ws.Cells.Item(22, 5) "  -0.50%  "
Set-TextValue $ws.Cells.Item(23, 2) "Dai"
Set-TextValue $ws.Cells.Item(23, 3) "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Cells.Item(23, 4) "1.000"
Set-TextValue $ws.Cells.Item(23, 5) "  +0.03%  "
Set-TextValue $ws.Cells.Item(24, 2) "Chainlink"
Set-TextValue $ws.Cells.Item(24, 3) "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Cells.Item(24, 4) "7.618"
Set-TextValue $ws.Cells.Item(24, 5) "  +3.11%  "
Set-TextValue $ws.Cells.Item(25, 2) "BinanceUSD"
Set-TextValue $ws.Cells.Item(25, 3) "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Cells.Item(25, 4) "1.000"
Set-TextValue $ws.Cells.Item(25, 5) "  -0.04%  "
Set-TextValue $ws.Cells.Item(26, 2) "Monero"
Set-TextValue $ws.Cells.Item(26, 3) "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Cells.Item(26, 4) "156.67"
Set-TextValue $ws.Cells.Item(26, 5) "  -1.02%  "
Set-TextValue $ws.Cells.Item(27, 2) "Stellar"
Set-TextValue $ws.Cells.Item(27, 3) "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Cells.Item(27, 4) "0.1397"
Set-TextValue $ws.Cells.Item(27, 5) "  -1.28%  "
Set-TextValue $ws.Cells.Item(28, 2) "Cosmos"
Set-TextValue $ws.Cells.Item(28, 3) "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Cells.Item(28, 4) "8.471"
Set-TextValue $ws.Cells.Item(28, 5) "  -0.11%  "
Set-TextValue $ws.Cells.Item(29, 2) "EthereumClassic"
Set-TextValue $ws.Cells.Item(29, 3) "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Cells.Item(29, 4) "17.77"
Set-TextValue $ws.Cells.Item(29, 5) "  -0.55%  "
Set-TextValue $ws.Cells.Item(30, 2) "PancakeSwap"
Set-TextValue $ws.Cells.Item(30, 3) "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Cells.Item(30, 4) "1.485"
Set-TextValue $ws.Cells.Item(30, 5) "  -0.51%  "
Set-TextValue $ws.Cells.Item(31, 2) "Hedera"
Set-TextValue $ws.Cells.Item(31, 3) "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Cells.Item(31, 4) "0.05923"
Set-TextValue $ws.Cells.Item(31, 5) "  -5.56%  "
Set-TextValue $ws.Cells.Item(32, 2) "Toncoin"
Set-TextValue $ws.Cells.Item(32, 3) "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Cells.Item(32, 4) "1.282"
Set-TextValue $ws.Cells.Item(32, 5) "  +0.70%  "
Set-TextValue $ws.Cells.Item(33, 2) "Filecoin"
Set-TextValue $ws.Cells.Item(33, 3) "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Cells.Item(33, 4) "4.132"
Set-TextValue $ws.Cells.Item(33, 5) "  -0.23%  "
Set-TextValue $ws.Cells.Item(34, 2) "InternetComputer(DFINITY)"
Set-TextValue $ws.Cells.Item(34, 3) "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Cells.Item(34, 4) "4.079"
Set-TextValue $ws.Cells.Item(34, 5) "  -0.15%  "
Set-TextValue $ws.Cells.Item(35, 2) "LidoDAOToken"
Set-TextValue $ws.Cells.Item(35, 3) "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Cells.Item(35, 4) "1.912"
Set-TextValue $ws.Cells.Item(35, 5) "  +3.08%  "
Set-TextValue $ws.Cells.Item(36, 2) "ARBITRUM"
Set-TextValue $ws.Cells.Item(36, 3) "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Cells.Item(36, 4) "1.178"
Set-TextValue $ws.Cells.Item(36, 5) "  +0.26%  "
Set-TextValue $ws.Cells.Item(37, 2) "ImmutableX"
Set-TextValue $ws.Cells.Item(37, 3) "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Cells.Item(37, 4) "0.7228"
Set-TextValue $ws.Cells.Item(37, 5) "  -0.67%  "
Set-TextValue $ws.Cells.Item(38, 2) "HuobiToken"
Set-TextValue $ws.Cells.Item(38, 3) "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Cells.Item(38, 4) "2.599"
Set-TextValue $ws.Cells.Item(38, 5) "  -0.28%  "
Set-TextValue $ws.Cells.Item(39, 2) "MXToken"
Set-TextValue $ws.Cells.Item(39, 3) "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Cells.Item(39, 4) "2.806"
Set-TextValue $ws.Cells.Item(39, 5) "  -1.33%  "
Set-TextValue $ws.Cells.Item(40, 2) "Maker"
Set-TextValue $ws.Cells.Item(40, 3) "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Cells.Item(40, 4) "1.219.55"
Set-TextValue $ws.Cells.Item(40, 5) "  +1.69%  "
Set-TextValue $ws.Cells.Item(41, 2) "VeChain"
Set-TextValue $ws.Cells.Item(41, 3) "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Cells.Item(41, 4) "0.01777"
Set-TextValue $ws.Cells.Item(41, 5) "  -0.12%  "
Set-TextValue $ws.Cells.Item(42, 2) "TrustWalletToken"
Set-TextValue $ws.Cells.Item(42, 3) "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Cells.Item(42, 4) "0.9127"
Set-TextValue $ws.Cells.Item(42, 5) "  -1.03%  "
Set-TextValue $ws.Cells.Item(43, 2) "FraxShare"
Set-TextValue $ws.Cells.Item(43, 3) "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Cells.Item(43, 4) "6.136"
Set-TextValue $ws.Cells.Item(43, 5) "  -0.25%  "
Set-TextValue $ws.Cells.Item(44, 2) "PaxDollar"
Set-TextValue $ws.Cells.Item(44, 3) "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Cells.Item(44, 4) "0.9999"
Set-TextValue $ws.Cells.Item(44, 5) "  -0.04%  "
Set-TextValue $ws.Cells.Item(45, 2) "RocketPoolETH"
Set-TextValue $ws.Cells.Item(45, 3) "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Cells.Item(45, 4) "2.002.96"
Set-TextValue $ws.Cells.Item(45, 5) "  -0.47%  "
Set-TextValue $ws.Cells.Item(46, 2) "Quant"
Set-TextValue $ws.Cells.Item(46, 3) "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Cells.Item(46, 4) "101.94"
Set-TextValue $ws.Cells.Item(46, 5) "  +0.09%  "
Set-TextValue $ws.Cells.Item(47, 2) "Aave"
Set-TextValue $ws.Cells.Item(47, 3) "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Cells.Item(47, 4) "67.30"
Set-TextValue $ws.Cells.Item(47, 5) "  +2.17%  "
Set-TextValue $ws.Cells.Item(48, 2) "Aptos"
Set-TextValue $ws.Cells.Item(48, 3) "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Cells.Item(48, 4) "7.367"
Set-TextValue $ws.Cells.Item(48, 5) "  +10.17%  "
Set-TextValue $ws.Cells.Item(49, 2) "TheSandbox"
Set-TextValue $ws.Cells.Item(49, 3) "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Cells.Item(49, 4) "0.4051"
Set-TextValue $ws.Cells.Item(49, 5) "  -0.04%  "
Set-TextValue $ws.Cells.Item(50, 2) "BabyDogeCoin"
Set-TextValue $ws.Cells.Item(50, 3) "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Cells.Item(50, 4) "0.00000000118"
Set-TextValue $ws.Cells.Item(50, 5) "  -1.46%  "
Set-TextValue $ws.Cells.Item(51, 2) "EnergySwap"
Set-TextValue $ws.Cells.Item(51, 3) "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Cells.Item(51, 4) "9.141"
Set-TextValue $ws.Cells.Item(51, 5) "  -0.45%  "
